$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.300.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.62%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.446.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.98%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.42%  "

$ws.Range("E9").Value = "  -1.67%  "

$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000173"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.893.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.389.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.455.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.98%  "

$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("E22").Value = "  -4.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "635.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.30%  "

$ws.Range("E27").Value = "  -0.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0955"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.98%  "

$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.40"
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.130"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.00%  "

$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("E35").Value = "  -5.49%  "

$ws.Range("E36").Value = "  -2.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "151.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.364"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.00%  "

$ws.Range("E41").Value = "  -3.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.95%  "

$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0305"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "151.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.14%  "

$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0501"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0901"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.49%  "
